$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns of interest: D=4 (Fecha), J=10 (Volumen), K=11 (Precio minimo),
# L=12 (Precio maximo), M=13 (Precio promedio ponderado), P=16 (Precio $/Kg).
# All other columns (A,B,C,E,F,G,H,I,N,O,Q,R) are identical on every data
# row in this sheet, so they don't need to move.

$firstRow = 27
$lastRow  = 120
$newLastRow = 121

# 1) Snapshot the current D/J/K/L/M/P values for every data row that will
#    be shifted, before anything is overwritten.
$dVals = @{}
$jVals = @{}
$kVals = @{}
$lVals = @{}
$mVals = @{}
$pVals = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $dVals[$r] = $ws.Cells.Item($r, 4).Value2
    $jVals[$r] = $ws.Cells.Item($r, 10).Value2
    $kVals[$r] = $ws.Cells.Item($r, 11).Value2
    $lVals[$r] = $ws.Cells.Item($r, 12).Value2
    $mVals[$r] = $ws.Cells.Item($r, 13).Value2
    $pVals[$r] = $ws.Cells.Item($r, 16).Value2
}

# 2) A brand new data point is inserted at row 27; every previously
#    existing row from 27..120 shifts down one row (to 28..121).
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 4).Value = $dVals[$src]
    $ws.Cells.Item($r, 10).Value = $jVals[$src]
    $ws.Cells.Item($r, 11).Value = $kVals[$src]
    $ws.Cells.Item($r, 12).Value = $lVals[$src]
    $ws.Cells.Item($r, 13).Value = $mVals[$src]
    $ws.Cells.Item($r, 16).Value = $pVals[$src]
}

# 3) Row 27 becomes the new data point.
$ws.Cells.Item($firstRow, 4).Value = 44620
$ws.Cells.Item($firstRow, 10).Value = 250
$ws.Cells.Item($firstRow, 11).Value = 23000
$ws.Cells.Item($firstRow, 12).Value = 23000
$ws.Cells.Item($firstRow, 13).Value = 23000
$ws.Cells.Item($firstRow, 16).Value = 920

# 4) A new row 121 is appended, carrying everything the old row 120 had
#    (it is a full copy of the constant columns too, since the sheet has
#    no formulas and every row shares the same schema/categorical values).
$ws.Cells.Item($newLastRow, 1).Value = $ws.Cells.Item($lastRow, 1).Value2
$ws.Cells.Item($newLastRow, 2).Value = $ws.Cells.Item($lastRow, 2).Value()
$ws.Cells.Item($newLastRow, 3).Value = $ws.Cells.Item($lastRow, 3).Value()
$ws.Cells.Item($newLastRow, 4).Value = $dVals[$lastRow]
$ws.Cells.Item($newLastRow, 4).NumberFormat = $ws.Cells.Item($lastRow, 4).NumberFormat()
$ws.Cells.Item($newLastRow, 5).Value = $ws.Cells.Item($lastRow, 5).Value2
$ws.Cells.Item($newLastRow, 6).Value = $ws.Cells.Item($lastRow, 6).Value2
$ws.Cells.Item($newLastRow, 7).Value = $ws.Cells.Item($lastRow, 7).Value()
$ws.Cells.Item($newLastRow, 8).Value = $ws.Cells.Item($lastRow, 8).Value()
$ws.Cells.Item($newLastRow, 9).Value = $ws.Cells.Item($lastRow, 9).Value()
$ws.Cells.Item($newLastRow, 10).Value = $jVals[$lastRow]
$ws.Cells.Item($newLastRow, 11).Value = $kVals[$lastRow]
$ws.Cells.Item($newLastRow, 12).Value = $lVals[$lastRow]
$ws.Cells.Item($newLastRow, 13).Value = $mVals[$lastRow]
$ws.Cells.Item($newLastRow, 14).Value = $ws.Cells.Item($lastRow, 14).Value()
$ws.Cells.Item($newLastRow, 15).Value = $ws.Cells.Item($lastRow, 15).Value()
$ws.Cells.Item($newLastRow, 16).Value = $pVals[$lastRow]
$ws.Cells.Item($newLastRow, 17).Value = $ws.Cells.Item($lastRow, 17).Value2
$ws.Cells.Item($newLastRow, 18).Value = $ws.Cells.Item($lastRow, 18).Value()
